$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Semana 9"
Write-Host ($wb.Worksheets.Count)
foreach ($ws in $wb.Worksheets) {
    Write-Host $ws.Name
}
